# Update the statsmodels OLS summary text blocks embedded in cell B2 of each
# of the 28 worksheets ("backwardElimination" results).
#
# The diff re-runs/re-saves the regression report the following day:
#   Date:  Wed, 01 Jan 2020  ->  Thu, 02 Jan 2020   (every sheet)
#   Time:  23:19:00          ->  20:48:52 (sheet 1) / 20:48:53 (sheets 2-28)
#
# Everything else in the report text (coefficients, AIC/BIC, etc.) is
# unchanged, so a targeted string replacement on the existing cell value is
# the safest way to reproduce the edit exactly.

$wb = $excel.ActiveWorkbook

$oldDate = "Date:                Wed, 01 Jan 2020"
$newDate = "Date:                Thu, 02 Jan 2020"
$oldTime = "Time:                        23:19:00"

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    if ($i -eq 1) {
        $newTime = "Time:                        20:48:52"
    } else {
        $newTime = "Time:                        20:48:53"
    }

    $cell = $ws.Range("B2")
    $text = $cell.Value()

    $text = $text -replace [regex]::Escape($oldDate), $newDate
    $text = $text -replace [regex]::Escape($oldTime), $newTime

    $cell.Value = $text

    # Re-pin the (wrap-text) row back to its original, max-out height so the
    # content refresh above doesn't leave a stray autofit side effect behind.
    $ws.Rows.Item(2).RowHeight = 409.5
}
